$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I2").Value = "CQLive-PENGUIN1SA4062-1.23"
$ws.Range("J2").Value = "CQLive-PENGUIN1SA4062-1.23.apk"

$ws.Range("H13").Select()
